# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to the refreshed coinranking.com snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.517.62"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "3.105.20"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.49%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.106.21"
$ws.Range("E8").Value = "  -5.02%  "
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("E10").Value = "  -7.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.27%  "
$ws.Range("E13").Value = "  -8.62%  "
$ws.Range("E14").Value = "  -10.05%  "
$ws.Range("D15").Value = "3.619.86"
$ws.Range("E15").Value = "  -4.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.115"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "63.563.40"
$ws.Range("E17").Value = "  -4.31%  "
$ws.Range("D18").Value = "3.109.60"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("E19").Value = "  -8.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.88%  "
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("E23").Value = "  -5.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.18%  "
$ws.Range("E29").Value = "  -12.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.80%  "
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  -6.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.73%  "
$ws.Range("E40").Value = "  -17.58%  "
$ws.Range("E41").Value = "  -8.78%  "
$ws.Range("E42").Value = "  -10.00%  "
$ws.Range("E43").Value = "  -6.16%  "
$ws.Range("D44").Value = "2.829.12"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -13.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.81%  "
$ws.Range("E50").Value = "  -5.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.12%  "
